# chore: fix docker ?
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")

# Row 67: already had date (45111), add type, hours, and work description
$ws.Range("B67").Value = "Implémentation"
$ws.Range("C67").Value = 8
$ws.Range("D67").Value = "Correction docker et test en production, mise à jour infos box (api)"

# Row 68: add date and type
$ws.Range("A68").Value = 45112
$ws.Range("B68").Value = "Implémentation"

# Update selection/active cell to D66 (as last user action before save)
$ws.Range("D66").Select()

$wb.Save()
